$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column A (test_number) changes from 1 to 2 for rows 2-11
$ws.Range("A2:A11").Value = 2

# Row 9 specific answer changes
$ws.Range("B9").Value = 6
$ws.Range("D9").Value = 4
$ws.Range("E9").Value = 3

# Update the active selection to match the author's final cursor position
$ws.Range("E9").Select()
